$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-10-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-07 Tuesday", 2)

# Update the answer table. Addressing cells directly by (row, column)
# avoids any ambiguity from values that coincidentally repeat elsewhere
# in the table (several old/new answers collide across rows).
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "32÷5=6, 2"
$t.Cell(1, 2).Range.Text  = "81÷2=40, 1"
$t.Cell(1, 3).Range.Text  = "51÷8=6, 3"
$t.Cell(1, 4).Range.Text  = "74÷4=18, 2"
$t.Cell(1, 5).Range.Text  = "69÷9=7, 6"

$t.Cell(5, 1).Range.Text  = "95÷5=19, 0"
$t.Cell(5, 2).Range.Text  = "26÷4=6, 2"
$t.Cell(5, 3).Range.Text  = "61÷8=7, 5"
$t.Cell(5, 4).Range.Text  = "88÷7=12, 4"
$t.Cell(5, 5).Range.Text  = "55÷3=18, 1"

$t.Cell(9, 1).Range.Text  = "94÷3=31, 1"
$t.Cell(9, 2).Range.Text  = "12÷8=1, 4"
$t.Cell(9, 3).Range.Text  = "26÷2=13, 0"
$t.Cell(9, 4).Range.Text  = "83÷8=10, 3"
$t.Cell(9, 5).Range.Text  = "67÷5=13, 2"

$t.Cell(13, 1).Range.Text = "89÷6=14, 5"
$t.Cell(13, 2).Range.Text = "20÷6=3, 2"
$t.Cell(13, 3).Range.Text = "96÷8=12, 0"
$t.Cell(13, 4).Range.Text = "29÷7=4, 1"
$t.Cell(13, 5).Range.Text = "16÷5=3, 1"

$t.Cell(17, 1).Range.Text = "57÷6=9, 3"
$t.Cell(17, 2).Range.Text = "37÷2=18, 1"
$t.Cell(17, 3).Range.Text = "86÷6=14, 2"
$t.Cell(17, 4).Range.Text = "15÷4=3, 3"
$t.Cell(17, 5).Range.Text = "13÷8=1, 5"
